# Insert one new daily price-record row for "Perejil" (Vega Modelo de Temuco)
# right above the existing row 367, pushing rows 367-431 down to 368-432
# (this also carries row 367's old number formatting down with it, matching
# Excel's default "insert copies format from the row above" behaviour).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(367).Insert()

$ws.Cells.Item(367, 1).Value  = 10
$ws.Cells.Item(367, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(367, 3).Value  = "La Araucanía"
$ws.Cells.Item(367, 4).Value  = 44951
$ws.Cells.Item(367, 5).Value  = 9
$ws.Cells.Item(367, 6).Value  = 100112044
$ws.Cells.Item(367, 7).Value  = "Perejil"
$ws.Cells.Item(367, 8).Value  = "Sin especificar"
$ws.Cells.Item(367, 9).Value  = "Primera"
$ws.Cells.Item(367, 10).Value = 35
$ws.Cells.Item(367, 11).Value = 5000
$ws.Cells.Item(367, 12).Value = 5000
$ws.Cells.Item(367, 13).Value = 5000
$ws.Cells.Item(367, 14).Value = "`$/docena de atados (3 kilos)"
$ws.Cells.Item(367, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(367, 16).Value = 1667
$ws.Cells.Item(367, 17).Value = 3
$ws.Cells.Item(367, 18).Value = "Hortaliza"
